# Auto-generated edit script applying the Cactuar_Profits.xlsx diff
# Updates derived (computed) columns H-N across several rows on each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 3334407
$ws.Range("I20").Value = 3334407
$ws.Range("K20").Value = 3334407
$ws.Range("M20").Value = -3334177
$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20936
$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20468
$ws.Range("H35").Value = 3334407
$ws.Range("I35").Value = 3334407
$ws.Range("K35").Value = 3334407
$ws.Range("M35").Value = -3334028
$ws.Range("H68").Value = 74500
$ws.Range("J68").Value = 74500
$ws.Range("L68").Value = 74500
$ws.Range("N68").Value = -75998
$ws.Range("H71").Value = 74500
$ws.Range("J71").Value = 74500
$ws.Range("L71").Value = 223500
$ws.Range("N71").Value = -230988
$ws.Range("H115").Value = 458
$ws.Range("I115").Value = 458
$ws.Range("K115").Value = 1374
$ws.Range("M115").Value = 193
$ws.Range("H127").Value = 2403.8235
$ws.Range("I127").Value = 1123.4
$ws.Range("J127").Value = 2937.3333
$ws.Range("K127").Value = 3370.2
$ws.Range("L127").Value = 8811.999899999999
$ws.Range("M127").Value = 1589.8
$ws.Range("N127").Value = -18731.9999
$ws.Range("H132").Value = 3658.8667
$ws.Range("I132").Value = 1310.3611
$ws.Range("K132").Value = 3931.0833
$ws.Range("M132").Value = -1401.0833
$ws.Range("H138").Value = 6151.7393
$ws.Range("J138").Value = 6249.591
$ws.Range("L138").Value = 18748.773
$ws.Range("N138").Value = -29028.773

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5008.35
$ws.Range("I61").Value = 4103.25
$ws.Range("K61").Value = 4103.25
$ws.Range("M61").Value = -3891.25
$ws.Range("H110").Value = 1086.6522
$ws.Range("I110").Value = 944.9
$ws.Range("K110").Value = 944.9
$ws.Range("M110").Value = 1100.1
$ws.Range("H132").Value = 4154.5093
$ws.Range("I132").Value = 1467.9143
$ws.Range("K132").Value = 4403.742899999999
$ws.Range("M132").Value = -1873.742899999999
$ws.Range("H136").Value = 5008.35
$ws.Range("I136").Value = 4103.25
$ws.Range("K136").Value = 12309.75
$ws.Range("M136").Value = -9759.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 320
$ws.Range("I64").Value = 100
$ws.Range("J64").Value = 356.66666
$ws.Range("K64").Value = 100
$ws.Range("L64").Value = 356.66666
$ws.Range("M64").Value = 125
$ws.Range("N64").Value = -806.66666
$ws.Range("H67").Value = 320
$ws.Range("I67").Value = 100
$ws.Range("J67").Value = 356.66666
$ws.Range("K67").Value = 100
$ws.Range("L67").Value = 356.66666
$ws.Range("M67").Value = 680
$ws.Range("N67").Value = -1916.66666
$ws.Range("H94").Value = 16667223
$ws.Range("I94").Value = 19231164
$ws.Range("K94").Value = 19231164
$ws.Range("M94").Value = -19230713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 849.2105
$ws.Range("I16").Value = 658
$ws.Range("K16").Value = 658
$ws.Range("M16").Value = -371
$ws.Range("H97").Value = 22160.25
$ws.Range("J97").Value = 22160.25
$ws.Range("L97").Value = 22160.25
$ws.Range("N97").Value = -24142.25
$ws.Range("H102").Value = 30241
$ws.Range("J102").Value = 30241
$ws.Range("L102").Value = 30241
$ws.Range("N102").Value = -35109
$ws.Range("H113").Value = 849.2105
$ws.Range("I113").Value = 658
$ws.Range("K113").Value = 658
$ws.Range("M113").Value = 1512
$ws.Range("H132").Value = 37045096
$ws.Range("I132").Value = 47623930
$ws.Range("J132").Value = 19182.75
$ws.Range("K132").Value = 142871790
$ws.Range("L132").Value = 57548.25
$ws.Range("M132").Value = -142869260
$ws.Range("N132").Value = -62608.25
$ws.Range("H141").Value = 110087.1
$ws.Range("I141").Value = 99999
$ws.Range("J141").Value = 110292.98
$ws.Range("K141").Value = 99999
$ws.Range("L141").Value = 110292.98
$ws.Range("M141").Value = -94819
$ws.Range("N141").Value = -120652.98

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 4774.8
$ws.Range("J52").Value = 4774.8
$ws.Range("L52").Value = 14324.4
$ws.Range("N52").Value = -14856.4
$ws.Range("H60").Value = 259.5
$ws.Range("I60").Value = 491
$ws.Range("J60").Value = 28
$ws.Range("K60").Value = 1473
$ws.Range("L60").Value = 84
$ws.Range("M60").Value = -1222
$ws.Range("N60").Value = -586
$ws.Range("H80").Value = 5999.2856
$ws.Range("J80").Value = 5999.2856
$ws.Range("L80").Value = 17997.8568
$ws.Range("N80").Value = -19869.8568
$ws.Range("H83").Value = 5999.2856
$ws.Range("J83").Value = 5999.2856
$ws.Range("L83").Value = 53993.5704
$ws.Range("N83").Value = -63353.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 32797.625
$ws.Range("J123").Value = 32797.625
$ws.Range("L123").Value = 32797.625
$ws.Range("N123").Value = -37697.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3797.21
$ws.Range("J7").Value = 4771.4
$ws.Range("L7").Value = 4771.4
$ws.Range("N7").Value = -4995.4
$ws.Range("H22").Value = 873.7
$ws.Range("I22").Value = 912.25
$ws.Range("J22").Value = 848
$ws.Range("K22").Value = 912.25
$ws.Range("L22").Value = 848
$ws.Range("M22").Value = -617.25
$ws.Range("N22").Value = -1438
$ws.Range("H27").Value = 873.7
$ws.Range("I27").Value = 912.25
$ws.Range("J27").Value = 848
$ws.Range("K27").Value = 912.25
$ws.Range("L27").Value = 848
$ws.Range("M27").Value = -805.25
$ws.Range("N27").Value = -1062
$ws.Range("H99").Value = 41129.5
$ws.Range("I99").Value = 7259
$ws.Range("J99").Value = 75000
$ws.Range("K99").Value = 7259
$ws.Range("L99").Value = 75000
$ws.Range("M99").Value = -4264
$ws.Range("N99").Value = -80990
$ws.Range("H100").Value = 1894.1111
$ws.Range("I100").Value = 1880.875
$ws.Range("K100").Value = 1880.875
$ws.Range("M100").Value = -1339.875
$ws.Range("H126").Value = 3797.21
$ws.Range("J126").Value = 4771.4
$ws.Range("L126").Value = 14314.2
$ws.Range("N126").Value = -19254.2
$ws.Range("H127").Value = 250056000
$ws.Range("J127").Value = 74665
$ws.Range("L127").Value = 74665
$ws.Range("N127").Value = -84585

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3933
$ws.Range("I81").Value = 4460.5
$ws.Range("J81").Value = 3581.3333
$ws.Range("K81").Value = 8921
$ws.Range("L81").Value = 7162.6666
$ws.Range("M81").Value = -7860
$ws.Range("N81").Value = -9284.6666
$ws.Range("H84").Value = 3933
$ws.Range("I84").Value = 4460.5
$ws.Range("J84").Value = 3581.3333
$ws.Range("K84").Value = 44605
$ws.Range("L84").Value = 35813.333
$ws.Range("M84").Value = -39301
$ws.Range("N84").Value = -46421.333
$ws.Range("H100").Value = 1306.3
$ws.Range("I100").Value = 868.5714
$ws.Range("J100").Value = 2327.6667
$ws.Range("K100").Value = 1737.1428
$ws.Range("L100").Value = 4655.3334
$ws.Range("M100").Value = -1196.1428
$ws.Range("N100").Value = -5737.3334
$ws.Range("H126").Value = 1419
$ws.Range("I126").Value = 1468.7142
$ws.Range("K126").Value = 4406.142599999999
$ws.Range("M126").Value = -1936.142599999999
